$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8806694746017456
$ws.Range("B1").Value = 1.619658350944519
$ws.Range("C1").Value = 6.061868667602539
$ws.Range("D1").Value = 1.861445188522339
$ws.Range("E1").Value = 1.122338652610779
